# Refactor auth-codes sheet: remove the "그룹" (group) column and blank out
# the per-member "인증 번호" (auth code) values, per the accompanying commit
# message ("remove group column ... remove auth codes from team table").
#
# Column layout before: A 팀 번호 | B 팀명 | C 그룹 | D 멤버 구분 | E LDAP 닉네임 | F 인증 번호
# Column layout after:  A 팀 번호 | B 팀명 | C 멤버 구분 | D LDAP 닉네임 | E 인증 번호 (blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop column C ("그룹") entirely; D/E/F shift left into C/D/E, carrying
# their widths (10/20/12) and values with them.
$ws.Range("C1").EntireColumn.Delete()

# Clear the auth-code values (now in column E) for every data row, leaving
# the "인증 번호" header in E1 untouched.
$ws.Range("E2:E103").Value = ""
